$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in new forecast values for column B, rows 110-121
$values = @{
    110 = 516.65019274929227
    111 = 459.25203286655994
    112 = 397.84397137970029
    113 = 398.16037878233846
    114 = 399.68264381252743
    115 = 399.25312666369229
    116 = 397.33526810467691
    117 = 397.00785842393844
    118 = 396.70599348008614
    119 = 397.67744701382151
    120 = 397.12623831080612
    121 = 396.19664348308919
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 2).Value = $values[$row]
}

# Update the sheet view's top-left cell and selection, matching the diff
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 105
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C110").Select()
